$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ACHR) updates
$ws.Range("K2").Value = 59.9
$ws.Range("N2").Value = 85.82376350509293

# Row 3 (JOBY) updates
$ws.Range("D3").Value = 14.43
$ws.Range("F3").Value = 10.66
$ws.Range("K3").Value = 56.9
$ws.Range("N3").Value = 85.82376350509293
